$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the demand-forecast date range for row 4 (이노그리드)
$ws.Range("B4").Value = "2024.04.22~04.26"

# Update the confirmed offering price for row 14 (신한스팩12호) from "-" to "2000"
$ws.Range("D14").Value = "2000"
